$wb = $excel.ActiveWorkbook

# --- "studies" sheet (sheet2): add a new "PMID" column in H ---
$wsStudies = $wb.Worksheets.Item("studies")
$wsStudies.Cells.Item(1, 8).Value = "PMID"

# --- "counts" sheet (sheet4): add a new "notes" column in F ---
$wsCounts = $wb.Worksheets.Item("counts")
$wsCounts.Cells.Item(1, 6).Value = "notes"

# Update selections to match the saved workbook state:
# "studies" keeps a leftover selection at H2 but is no longer the active tab.
$wsStudies.Select()
$wsStudies.Range("H2").Select()

# "counts" becomes the active/visible tab, selection resting on F2.
$wsCounts.Select()
$wsCounts.Range("F2").Select()
